$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 0.30952380952380953
$ws.Range("G2").Value = 13
$ws.Range("H2").Value = 0.40476190476190477
$ws.Range("I2").Value = 17
$ws.Range("D5").Value = 0.043478260869565216
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.086956521739130432
$ws.Range("G5").Value = 4
$ws.Range("N5").Value = 0.2391304347826087
$ws.Range("O5").Value = 11
$ws.Range("J6").Value = 0.046511627906976744
$ws.Range("K6").Value = 2
$ws.Range("D7").Value = 0.12
$ws.Range("E7").Value = 3
$ws.Range("J8").Value = 0.054545454545454543
$ws.Range("K8").Value = 3
$ws.Range("N8").Value = 0.21818181818181817
$ws.Range("O8").Value = 12
$ws.Range("F9").Value = 0.0625
$ws.Range("G9").Value = 2
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("D10").Value = 0.026315789473684209
$ws.Range("E10").Value = 1
$ws.Range("L10").Value = 0.15789473684210525
$ws.Range("M10").Value = 6
$ws.Range("D11").Value = 0.078947368421052627
$ws.Range("E11").Value = 3
$ws.Range("D12").Value = 0.28205128205128205
$ws.Range("E12").Value = 11
$ws.Range("C15").Value = 23
$ws.Range("D15").Value = 0.086956521739130432
$ws.Range("F15").Value = 0.17391304347826086
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 0.2608695652173913
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 0.043478260869565216
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.043478260869565216
$ws.Range("M15").Value = 1
$ws.Range("N15").Value = 0.086956521739130432
$ws.Range("O15").Value = 2
$ws.Range("H17").Value = 0.375
$ws.Range("I17").Value = 6
$ws.Range("D18").Value = 0.23076923076923078
$ws.Range("E18").Value = 6
$ws.Range("F18").Value = 0.34615384615384615
$ws.Range("G18").Value = 9
$ws.Range("D19").Value = 0.28846153846153844
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = 0.32692307692307693
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 0.38461538461538464
$ws.Range("I19").Value = 20
$ws.Range("N19").Value = 0.19230769230769232
$ws.Range("O19").Value = 10
$ws.Range("L20").Value = 0.057692307692307696
$ws.Range("M20").Value = 3
$ws.Range("N20").Value = 0.15384615384615385
$ws.Range("O20").Value = 8
$ws.Range("H21").Value = 0.42857142857142855
$ws.Range("I21").Value = 18
$ws.Range("L22").Value = 0.066666666666666666
$ws.Range("M22").Value = 3
$ws.Range("N23").Value = 0.095238095238095233
$ws.Range("O23").Value = 4
$ws.Range("L24").Value = 0.05128205128205128
$ws.Range("M24").Value = 2
$ws.Range("F26").Value = 0.068965517241379309
$ws.Range("G26").Value = 2
$ws.Range("N26").Value = 0.2413793103448276
$ws.Range("O26").Value = 7
$ws.Range("L28").Value = 0.25
$ws.Range("M28").Value = 13
$ws.Range("D29").Value = 0.088888888888888892
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 0.17777777777777778
$ws.Range("G29").Value = 8
$ws.Range("L29").Value = 0.22222222222222221
$ws.Range("M29").Value = 10
$ws.Range("D30").Value = 0.14285714285714285
$ws.Range("E30").Value = 6
$ws.Range("J30").Value = 0.19047619047619047
$ws.Range("K30").Value = 8
$ws.Range("D31").Value = 0.17647058823529413
$ws.Range("E31").Value = 6
$ws.Range("F31").Value = 0.26470588235294118
$ws.Range("G31").Value = 9
$ws.Range("J31").Value = 0.23529411764705882
$ws.Range("K31").Value = 8
$ws.Range("L31").Value = 0.26470588235294118
$ws.Range("M31").Value = 9
$ws.Range("N31").Value = 0.44117647058823528
$ws.Range("O31").Value = 15
$ws.Range("F32").Value = 0.063492063492063489
$ws.Range("G32").Value = 4
$ws.Range("H32").Value = 0.14285714285714285
$ws.Range("I32").Value = 9
$ws.Range("H34").Value = 0.12820512820512819
$ws.Range("I34").Value = 5
$ws.Range("N34").Value = 0.076923076923076927
$ws.Range("O34").Value = 3
$ws.Range("D36").Value = 0.1
$ws.Range("E36").Value = 2
$ws.Range("J37").Value = 0.28000000000000003
$ws.Range("K37").Value = 7
$ws.Range("D38").Value = 0.0625
$ws.Range("E38").Value = 1
$ws.Range("H40").Value = 0.59375
$ws.Range("I40").Value = 19
$ws.Range("L40").Value = 0.21875
$ws.Range("M40").Value = 7
$ws.Range("N40").Value = 0.25
$ws.Range("O40").Value = 8
$ws.Range("N41").Value = 0.071428571428571425
$ws.Range("O41").Value = 3
$ws.Range("D42").Value = 0.021276595744680851
$ws.Range("E42").Value = 1
$ws.Range("J42").Value = 0.042553191489361701
$ws.Range("K42").Value = 2
$ws.Range("J43").Value = 0.10416666666666667
$ws.Range("K43").Value = 10
$ws.Range("L43").Value = 0.20833333333333334
$ws.Range("M43").Value = 20
$ws.Range("D44").Value = 0.12
$ws.Range("E44").Value = 3
$ws.Range("H44").Value = 0.2
$ws.Range("I44").Value = 5
$ws.Range("L44").Value = 0.08
$ws.Range("M44").Value = 2
$ws.Range("F45").Value = 0.31034482758620691
$ws.Range("G45").Value = 9
$ws.Range("F46").Value = 0.2608695652173913
$ws.Range("G46").Value = 6
$ws.Range("F47").Value = 0.060606060606060608
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 0.12121212121212122
$ws.Range("I47").Value = 4
$ws.Range("N48").Value = 0.17391304347826086
$ws.Range("O48").Value = 8
$ws.Range("H49").Value = 0.071428571428571425
$ws.Range("I49").Value = 4
$ws.Range("N49").Value = 0.089285714285714288
$ws.Range("O49").Value = 5
$ws.Range("N50").Value = 0.19148936170212766
$ws.Range("O50").Value = 9
$ws.Range("D51").Value = 0.13333333333333333
$ws.Range("E51").Value = 6
$ws.Range("F51").Value = 0.17777777777777778
$ws.Range("G51").Value = 8
$ws.Range("H51").Value = 0.17777777777777778
$ws.Range("I51").Value = 8
$ws.Range("J51").Value = 0.088888888888888892
$ws.Range("K51").Value = 4
$ws.Range("N51").Value = 0.24444444444444444
$ws.Range("O51").Value = 11
$ws.Range("J52").Value = 0.05
$ws.Range("K52").Value = 2
$ws.Range("F53").Value = 0.11764705882352941
$ws.Range("G53").Value = 4
$ws.Range("L53").Value = 0.058823529411764705
$ws.Range("M53").Value = 2
$ws.Range("F54").Value = 0.13157894736842105
$ws.Range("G54").Value = 5
$ws.Range("N54").Value = 0.078947368421052627
$ws.Range("O54").Value = 3
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("N55").Value = 0.31034482758620691
$ws.Range("O55").Value = 9
